$d = $word.ActiveDocument

# 1. "College Database" -> "CSUF College Database"
$d.Content.Find.Execute(
    "College Database", $true, $false, $false, $false, $false,
    $true, 1, $false, "CSUF College Database", 2) | Out-Null

# 2. College Database tech line: "Vite, React, Bootstrap, Laravel" -> "React/Vite, Bootstrap, Laravel"
$d.Content.Find.Execute(
    "Vite, React, Bootstrap, Laravel", $true, $false, $false, $false, $false,
    $true, 1, $false, "React/Vite, Bootstrap, Laravel", 2) | Out-Null

# 3. "Web database application for a university" -> "Web database application for the university"
$d.Content.Find.Execute(
    "Web database application for a university", $true, $false, $false, $false, $false,
    $true, 1, $false, "Web database application for the university", 2) | Out-Null

# 4. "Web Development:" bullet tech list: "Vite, React, Next.js" -> "React/Vite/Next.js"
$d.Content.Find.Execute(
    "Vite, React, Next.js", $true, $false, $false, $false, $false,
    $true, 1, $false, "React/Vite/Next.js", 2) | Out-Null

# 5. "OpenAI, Stripe, Clerk, " stays the same text but two runs merge into one (no visible text change)
$d.Content.Find.Execute(
    "OpenAI, Stripe, Clerk, ", $true, $false, $false, $false, $false,
    $true, 1, $false, "OpenAI, Stripe, Clerk, ", 2) | Out-Null
